$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers): split the single property_id/new_value pair into two pairs ---
$ws.Range("B1").Value = "property_id1"
$ws.Range("C1").Value = "new_value1"
$ws.Range("D1").Value = "property_id2"
$ws.Range("E1").Value = "new_value2"

# --- Row 2 (data): update accession number and first field, add second field ---
$ws.Range("A2").Value = 202201037
$ws.Range("B2").Value = "catalogBlock_Fish_GenderId"
$ws.Range("C2").Value = "Hona"
$ws.Range("D2").Value = "catalogBlock_Fish_Totallength"
$ws.Range("E2").Value = 123

# --- New column widths for the two new columns D and E ---
$ws.Columns.Item(4).ColumnWidth = 30.5
$ws.Columns.Item(5).ColumnWidth = 31

# --- Selection moves to A2 ---
$ws.Range("A2").Select()
